$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "57.840.33"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +2.08%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.334.14"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +0.67%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "538.66"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +3.15%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "134.05"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +1.69%  "
$ws.Range("E7").Value = "  +0.42%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.556"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +4.24%  "
$ws.Range("E9").Value = "  -0.02%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "5.52"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +4.10%  "
$ws.Range("E11").Value = "  -0.26%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.355"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +4.22%  "
$ws.Range("E13").Value = "  +0.60%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.745.80"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +0.13%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "57.761.79"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +1.94%  "
$ws.Range("E16").Value = "  +0.37%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.337.41"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +0.55%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "10.65"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +2.14%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "331.74"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -1.26%  "
$ws.Range("E20").Value = "  +2.05%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.67"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -2.34%  "
$ws.Range("E22").Value = "  -0.03%  "
$ws.Range("E23").Value = "  +0.16%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "62.58"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +1.48%  "
$ws.Range("E25").Value = "  +1.89%  "
$ws.Range("E26").Value = "  -2.39%  "
$ws.Range("E27").Value = "  +0.72%  "
$ws.Range("E28").Value = "  +7.67%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.76"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +3.63%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "170.48"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +0.70%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0₃0731"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +1.44%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.08"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -0.55%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "18.47"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +0.43%  "
$ws.Range("E34").Value = "  +12.70%  "
$ws.Range("E35").Value = "  +0.01%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.998"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +0.44%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.21"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +5.75%  "
$ws.Range("E38").Value = "  -0.93%  "
$ws.Range("E39").Value = "  +3.74%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "38.87"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +0.05%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "144.47"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -3.19%  "
$ws.Range("E42").Value = "  -0.01%  "
$ws.Range("E43").Value = "  +1.03%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "285.11"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -0.39%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0939"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +1.10%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "19.01"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +2.43%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.561"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +0.44%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.385"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +1.51%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0217"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +0.73%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "17.47"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +1.48%  "
